$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "GRT-USD"
$ws.Range("A17").Value = "BSCX-USD"
